$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 12: center-align the existing cells (A, B, D now horizontal+vertical
#    center like the other row-header cells; E gets horizontal center added
#    on top of its existing wrap; C gets horizontal center too).
# ---------------------------------------------------------------------------
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("C12").HorizontalAlignment = -4108
$ws.Range("D12").HorizontalAlignment = -4108
$ws.Range("E12").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 2. New rows 13-17 continue the "marco-v1-llama3.2:3b" conversation turns in
#    columns F (user query) / G (response); columns A-E stay blank but carry
#    the same centered formatting as the rest of the merged block.
# ---------------------------------------------------------------------------
$ws.Range("F13").Value = 'Italian'
$ws.Range("G13").Value = '**Delicious!** Just Italian. And what kind of Italian dish are you thinking of? (e.g., pasta, pizza, risotto)'

$ws.Range("F14").Value = 'I want to eat pizza'
$ws.Range("G14").Value = '**Great Choice!** Pizza! And what kind of pizza are you in the mood for? (e.g., Margherita, Pepperoni, Veggie) Or
maybe something more unique?'

$ws.Range("F15").Value = 'Veggie'
$ws.Range("G15").Value = '**Veggie Pizza!** And how many people are you looking to feed? (e.g., just me, my partner, a big group?)'

$ws.Range("F16").Value = 'for 4 people'
$ws.Range("G16").Value = '**For 4 People!** Okay, so a veggie pizza for 4 people. And what is your budget like for this meal? (e.g., 10-20,
20-50, 50+?)'

$ws.Range("F17").Value = 200
$ws.Range("G17").Value = '**Within Your Budget!** Okay, so a veggie pizza for 4 people within a budget of 200. I''ll find the perfect place
for you!'

# Column A-E formatting for the new rows (matches the rest of the block).
$ws.Range("A13:A17").HorizontalAlignment = -4108
$ws.Range("B13:B17").HorizontalAlignment = -4108
$ws.Range("C13:C17").HorizontalAlignment = -4108
$ws.Range("D13:D17").HorizontalAlignment = -4108
$ws.Range("E13:E17").HorizontalAlignment = -4108

# Column G formatting for the wrapped response cells (row 15 stays plain).
$ws.Range("G13").WrapText = $true
$ws.Range("G13").VerticalAlignment = -4108
$ws.Range("G14").WrapText = $true
$ws.Range("G14").VerticalAlignment = -4108
$ws.Range("G16").WrapText = $true
$ws.Range("G16").VerticalAlignment = -4108
$ws.Range("G17").WrapText = $true
$ws.Range("G17").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Row heights for the new rows (row 15 keeps the sheet default).
# ---------------------------------------------------------------------------
$ws.Rows(13).RowHeight = 45
$ws.Rows(14).RowHeight = 60
$ws.Rows(16).RowHeight = 60
$ws.Rows(17).RowHeight = 60

# ---------------------------------------------------------------------------
# 4. Merge the new row span (12:17) for columns A-E, matching the existing
#    merged blocks used for rows 2-6 and 7-11.
# ---------------------------------------------------------------------------
$ws.Range("A12:A17").Merge()
$ws.Range("B12:B17").Merge()
$ws.Range("C12:C17").Merge()
$ws.Range("D12:D17").Merge()
$ws.Range("E12:E17").Merge()

# ---------------------------------------------------------------------------
# 5. Scroll / selection so the new rows are in view, matching the author's
#    saved viewport.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D12:D17").Select()
